$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Indie Project" log entry to mention the additional work on
# insert/delete cascade testing, and log 4 hours for that day.
$ws.Range("D32").Value = "Indie Project: figuring out getByProperty issue in StoryDaoTest; reviewing generic dao examples; worked on testing effects of insert and delete on associated entities`nWeek 5: watched follow up videos"
$ws.Range("B32").Value = 4
$ws.Range("A32:D32").RowHeight = 45

# Select the two short, now-redundant notes ("1+ hr 1st thing" / "7:05-x")
# and delete the entire rows, shifting everything below up by two.
$ws.Range("A34:XFD35").Select()
$ws.Range("A34:D35").EntireRow.Delete()
